# Update automatico del mapa (AYKO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix E80: was an empty inline string, now has OT value ---
$ws.Range("E80").Value = "Pendiente ADM"

# --- Append new rows 84-86 ---
# Columns A-H, J-L, O-P are text; I, M, N are numeric.
$textCols = @("A","B","C","D","E","F","G","H","J","K","L","O","P")

$newRows = @(
    @{ Row = 84; A = "6469"; B = "8/6/2025"; C = "USPALLATA 3626"; D = "4"; E = "808733910"; F = "AYKO"; G = "Pendiente"; H = "Aplomar"; I = 1; J = "Aplomo"; K = "Sin equipos"; L = "Pasante"; M = -58.414389; N = -34.641308; O = "Boedo"; P = "Capital Sur" },
    @{ Row = 85; A = "6561"; B = "8/6/2025"; C = "GOLETA SARANDI 6050 "; D = "8"; E = "808733912"; F = "AYKO"; G = "Pendiente"; H = "Picada"; I = 1; J = "Cambio"; K = "Sin equipos"; L = "Pasante"; M = -58.468841; N = -34.686635; O = "Boedo"; P = "Capital Sur" },
    @{ Row = 86; A = "6567"; B = "8/6/2025"; C = "SAN ANTONIO 1221"; D = "4"; E = "808733914"; F = "AYKO"; G = "Pendiente"; H = "Aplomar"; I = 1; J = "Aplomo"; K = "Sin equipos"; L = "Pasante"; M = -58.375684; N = -34.656092; O = "San Telmo"; P = "Capital Sur" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Force the text columns to be stored as text (not auto-converted to
    # numbers/dates) by temporarily applying a text number format to just
    # those cells, then restore the default "Normal" style so no stray
    # formatting remains. Numeric columns (I, M, N) are left untouched so
    # they are written as real numbers.
    $textRange = $ws.Range("A" + $rowNum + ":H" + $rowNum + "," + "J" + $rowNum + ":L" + $rowNum + "," + "O" + $rowNum + ":P" + $rowNum)
    $textRange.NumberFormat = "@"

    foreach ($col in $textCols) {
        $ws.Range($col + $rowNum).Value = [string]$r[$col]
    }

    $ws.Range("I" + $rowNum).Value = $r.I
    $ws.Range("M" + $rowNum).Value = $r.M
    $ws.Range("N" + $rowNum).Value = $r.N

    $textRange.Style = "Normal"
}
